# PAS-2715: Updated VIN Upload files
#  - Remove the STAT column (old col Z) and the CHOICE_TIER column (old col AC)
#  - Append eight new columns after the (now shifted) ALTFUEL column:
#       BI_SYMBOL, PD_SYMBOL, UM_SYMBOL, MP_SYMBOL, ENTRYDATE, VALID,
#       ANTITHEFT_DISCOUNT, RESTRAINTS_DISCOUNT
#  - Populate row 2 sample data for the new columns
#  - Update the view: scroll right, select AK11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two retired columns -----------------------------------
# Deleting entire columns shifts everything to the right of them left by
# one, which is exactly what the diff shows happened to COLL_SYMBOL,
# COMP_SYMBOL and ALTFUEL (and their row-2 values).
$ws.Range("Z1").EntireColumn.Delete()    # old STAT column
$ws.Range("AB1").EntireColumn.Delete()   # old CHOICE_TIER column (now at AB after first delete)

# --- Add the new trailing columns --------------------------------------
# After the two deletions the data now ends at column AB (ALTFUEL).
# New columns go in AC:AJ.
$ws.Range("AC1").Value2 = "BI_SYMBOL"
$ws.Range("AD1").Value2 = "PD_SYMBOL"
$ws.Range("AE1").Value2 = "UM_SYMBOL"
$ws.Range("AF1").Value2 = "MP_SYMBOL"
$ws.Range("AG1").Value2 = "ENTRYDATE"
$ws.Range("AH1").Value2 = "VALID"
$ws.Range("AI1").Value2 = "ANTITHEFT_DISCOUNT"
$ws.Range("AJ1").Value2 = "RESTRAINTS_DISCOUNT"

# Match the header style ("Good" -- green fill) used by the rest of row 1.
$ws.Range("AC1:AJ1").Style = "Good"

# Row 2 sample values for the new columns.
$ws.Range("AC2").Value2 = "K"
$ws.Range("AD2").Value2 = "K"
$ws.Range("AE2").Value2 = "K"
$ws.Range("AF2").Value2 = "K"
$ws.Range("AG2").Value2 = 20000101
$ws.Range("AH2").Value2 = "Y"
$ws.Range("AI2").Value2 = "Y"
$ws.Range("AJ2").Value2 = "N"

# These cells share the left-aligned style ("s=3") used throughout row 2
# (everything except ENTRYDATE, which stays unstyled/general).
$ws.Range("AC2:AF2").HorizontalAlignment = -4131   # xlLeft
$ws.Range("AH2:AJ2").HorizontalAlignment = -4131   # xlLeft

# --- View state ----------------------------------------------------------
$ws.Range("AK11").Select()
$excel.ActiveWindow.ScrollColumn = 25
